$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# remember original row heights so autofit from the multi-line text
# assignments below does not change the saved row heights
$row14Height = $ws.Rows.Item(14).RowHeight
$row17Height = $ws.Rows.Item(17).RowHeight

$ws.Range("D11").Value = "5. Verifica Dados"
$ws.Range("D12").Value = "6. Regista Cliente"
$ws.Range("D13").Value = "7. Informa cliente de que foi registado com sucesso"
$ws.Range("B14").Value = "Alternativa 1`n[Cliente já registado]`n(Passo 6)"
$ws.Range("D14").Value = "6.1. Informa cliente que já está registado"
$ws.Range("D15").Value = "6.2. Sai do ecrã de registo"
$ws.Range("D16").Value = "Regressa a 1"
$ws.Range("B17").Value = "Alternativa 2`n[Dados inválidos]`n(Passo 6)"
$ws.Range("D17").Value = "6.1.1. Informa cliente que os dados são inválidos"
$ws.Range("D18").Value = "Regressa a 4"

$ws.Rows.Item(14).RowHeight = $row14Height
$ws.Rows.Item(17).RowHeight = $row17Height

$ws.Range("D17").Select()
